$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Prix Spot": add a new column BK (63) = "15-aug" with the
# 24 hourly spot prices beneath it.
# -----------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Header cell, formatted like the rest of row 1 (bold / centered / bordered)
$wsSpot.Cells.Item(1, 63).Value = "15-aug"
$wsSpot.Range("BJ1").Copy()
$wsSpot.Cells.Item(1, 63).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$spotValues = @(
    90.52,
    83.53,
    79.41,
    68.47,
    76.09,
    65.65000000000001,
    73.45999999999999,
    79.98999999999999,
    71.89,
    43.97,
    25,
    3,
    0.08,
    0,
    0,
    1.72,
    23.53,
    33.23,
    61.12,
    100.82,
    111.39,
    102,
    102.89,
    93.45
)

$row = 2
foreach ($v in $spotValues) {
    $wsSpot.Cells.Item($row, 63).Value = $v
    $row = $row + 1
}

# -----------------------------------------------------------------
# Sheet "Gaz": append row 60 with the new daily price.
# The date is written as plain text (like the existing rows) rather
# than letting Excel auto-convert it to a date serial number.
# -----------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A60").Formula = '="2025-08-13"'
$wsGaz.Range("A60").Copy()
$wsGaz.Range("A60").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$wsGaz.Range("B60").Value = 31.85

# -----------------------------------------------------------------
# Sheet "CO2": append row 60 with the new daily price.
# -----------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A60").Formula = '="2025-08-13"'
$wsCO2.Range("A60").Copy()
$wsCO2.Range("A60").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$wsCO2.Range("B60").Value = 71.06999999999999
